$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = 27784
$ws.Range("F4").Value = 28447
$ws.Range("G4").Value = 29070
$ws.Range("H4").Value = 29523

$ws.Range("E4:H4").Select()
$excel.ActiveCell = $ws.Range("E4")
